$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows documenting wiring worked on over the summer (Automation Direct cable
# + silicone hookup wire), inserted right after the existing "Wiring" section entries.
$ws.Range("A5").Value = "Automation Direct"
$ws.Range("D5").Value = "22 AWG 4-Core Signal Cable Unshielded"
$ws.Range("D6").Value = "Silicone wire red"
$ws.Range("D7").Value = "Silicone wire black"

# Normalize "N pin" -> "N-pin" wording in the Molex Micro-Fit descriptions, and
# clarify the right-angle connectors as PCB-mount.
$ws.Range("D12").Value = "3-pin Molex Micro-Fit receptacle connector"
$ws.Range("D13").Value = "3-pin Molex Micro-Fit plug connector"
$ws.Range("D14").Value = "4-pin Molex Micro-Fit receptacle connector"
$ws.Range("D15").Value = "4-pin Molex Micro-Fit plug connector"
$ws.Range("D16").Value = "6-pin Molex Micro-Fit receptacle connector"
$ws.Range("D17").Value = "6-pin Molex Micro-Fit plug connector"
$ws.Range("D18").Value = "3-pin PCB-mount right angle Molex microfit receptacle connector"
$ws.Range("D19").Value = "4-pin PCB-mount right angle Molex microfit receptacle connector"
$ws.Range("D20").Value = "6-pin PCB-mount right angle Molex microfit receptacle connector"

# Restore the view to the top of the sheet with the frozen header row, matching
# where the author left the selection after this edit.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("A42").Select()
